$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.048.37'
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = '1.763.22'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '237.81'
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").Value = '0.5228'
$ws.Range("E7").Value = '  +2.75%  '
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("D9").Value = '40.45'
$ws.Range("E9").Value = '  -3.91%  '
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").Value = '1.775.06'
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("D12").Value = '0.07023'
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").Value = '15.72'
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("E14").Value = '  +6.20%  '
$ws.Range("D15").Value = '4.546'
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").Value = '78.21'
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("D19").Value = '26.072.03'
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").Value = '0.000006734'
$ws.Range("E21").Value = '  -3.19%  '
$ws.Range("D22").Value = '2.002.61'
$ws.Range("E22").Value = '  +1.64%  '
$ws.Range("D23").Value = '4.077'
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("E24").Value = '  +2.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.200'
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").Value = '138.77'
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("D27").Value = '1.517'
$ws.Range("E27").Value = '  +3.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.840'
$ws.Range("E28").Value = '  +0.98%  '
$ws.Range("D29").Value = '15.19'
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("D30").Value = '103.46'
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("D31").Value = '0.08393'
$ws.Range("E31").Value = '  +2.80%  '
$ws.Range("D32").Value = '3.702'
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("D33").Value = '3.454'
$ws.Range("E33").Value = '  -1.69%  '
$ws.Range("E34").Value = '  -1.50%  '
$ws.Range("D35").Value = '2.624'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = '1.003'
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("D37").Value = '0.6063'
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("D38").Value = '2.742'
$ws.Range("E38").Value = '  +2.13%  '
$ws.Range("D39").Value = '0.01592'
$ws.Range("E39").Value = '  +2.27%  '
$ws.Range("D40").Value = '1.989'
$ws.Range("E40").Value = '  +3.87%  '
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("D42").Value = '102.78'
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("D43").Value = '0.3882'
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.940'
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").Value = '0.05515'
$ws.Range("E46").Value = '  +2.18%  '
$ws.Range("D47").Value = '6.345'
$ws.Range("E47").Value = '  +5.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1120'
$ws.Range("E48").Value = '  +0.62%  '
$ws.Range("D49").Value = '30.23'
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").Value = '52.71'
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("D51").Value = '1.005'
$ws.Range("E51").Value = '  +0.93%  '
